$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new "Save" column (copy formatting from the adjacent "sum" header)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-12
$saveValues = @(0, 0, 0, 1, 0, 1, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
